$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'68.887.18"
$ws.Range("E2").Value = '  +1.49%  '
$ws.Range("D3").Value = "'3.357.99"
$ws.Range("E3").Value = '  +1.12%  '
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = "'585.17"
$ws.Range("E5").Value = '  +1.03%  '
$ws.Range("D6").Value = "'178.97"
$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = '  -0.10%  '
$ws.Range("D8").Value = "'0.593"
$ws.Range("E8").Value = '  +0.95%  '
$ws.Range("E9").Value = '  +4.50%  '
$ws.Range("D10").Value = "'0.584"
$ws.Range("E10").Value = '  +1.56%  '
$ws.Range("D11").Value = "'48.17"
$ws.Range("E11").Value = '  +6.02%  '
$ws.Range("E12").Value = '  +2.27%  '
$ws.Range("D13").Value = "'694.80"
$ws.Range("E13").Value = '  +5.60%  '
$ws.Range("D14").Value = "'3.914.61"
$ws.Range("E14").Value = '  +1.30%  '
$ws.Range("D15").Value = "'8.49"
$ws.Range("E15").Value = '  +1.14%  '
$ws.Range("D16").Value = "'68.897.09"
$ws.Range("E16").Value = '  +1.57%  '
$ws.Range("E17").Value = '  +1.42%  '
$ws.Range("D18").Value = "'3.365.66"
$ws.Range("E18").Value = '  +1.39%  '
$ws.Range("D19").Value = "'17.57"
$ws.Range("E19").Value = '  +1.22%  '
$ws.Range("E20").Value = '  +2.85%  '
$ws.Range("D21").Value = "'0.897"
$ws.Range("E21").Value = '  +1.38%  '
$ws.Range("D22").Value = "'5.47"
$ws.Range("E22").Value = '  +2.52%  '
$ws.Range("D23").Value = "'17.06"
$ws.Range("E23").Value = '  +0.99%  '
$ws.Range("D24").Value = "'101.34"
$ws.Range("E24").Value = '  +3.36%  '
$ws.Range("E25").Value = '  +1.74%  '
$ws.Range("E26").Value = '  +1.61%  '
$ws.Range("D27").Value = "'9.54"
$ws.Range("E27").Value = '  +3.15%  '
$ws.Range("D28").Value = "'33.48"
$ws.Range("E28").Value = '  +0.99%  '
$ws.Range("D29").Value = "'8.59"
$ws.Range("E29").Value = '  +2.17%  '
$ws.Range("D30").Value = "'7.03"
$ws.Range("E30").Value = '  -2.60%  '
$ws.Range("D31").Value = "'11.11"
$ws.Range("E31").Value = '  +1.83%  '
$ws.Range("D32").Value = "'553.66"
$ws.Range("E32").Value = '  -2.60%  '
$ws.Range("E33").Value = '  +1.23%  '
$ws.Range("D34").Value = "'3.57"
$ws.Range("E34").Value = '  +10.69%  '
$ws.Range("D35").Value = "'57.73"
$ws.Range("E35").Value = '  +2.44%  '
$ws.Range("D37").Value = "'3.709.33"
$ws.Range("E37").Value = '  +0.77%  '
$ws.Range("E38").Value = '  +8.57%  '
$ws.Range("E39").Value = '  +1.54%  '
$ws.Range("E40").Value = '  +3.55%  '
$ws.Range("D41").Value = "'2.63"
$ws.Range("E41").Value = '  +0.60%  '
$ws.Range("D42").Value = "'0.0₃0678"
$ws.Range("E42").Value = '  +2.75%  '
$ws.Range("D43").Value = "'0.338"
$ws.Range("E43").Value = '  +1.87%  '
$ws.Range("D44").Value = "'0.0415"
$ws.Range("E44").Value = '  +2.70%  '
$ws.Range("D45").Value = "'3.26"
$ws.Range("E45").Value = '  -1.47%  '
$ws.Range("D46").Value = "'2.65"
$ws.Range("E46").Value = '  +2.38%  '
$ws.Range("E47").Value = '  +1.24%  '
$ws.Range("E48").Value = '  +0.04%  '
$ws.Range("E49").Value = '  -1.76%  '
$ws.Range("D50").Value = "'131.84"
$ws.Range("E50").Value = '  +3.91%  '
$ws.Range("E51").Value = '  -1.94%  '
